# Updates to support ORA compression and the move to the hci-dragen server.
#
# - The D-column "Extra MultiQC Options" config path cells drop the
#   "Test/" path segment (server move off the test box).
# - The 10x Genomics Library Kit cell (B10/B11) gains a new kit entry.
# - View state: freeze-pane top-left cell and active selection move.
# - Row 4 header grows taller; rows 10/11 grow taller to fit the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newConfigPath = "--config /home/tomatosrvs/AutoAnalysis/autoAnalysis.multiqc.config.yaml"

$new10xKit = "10X Genomics Next GEM Single Cell  3' Gene Expression Library prep v3.1  with UDI; 10X Genomics Sigle Cell 3' Cell Multiplexing with UDI; 10x Genomics Chromium Singel Cell Fixed RNA Profiling; 10x Genomics Chromium Next GEM Single Cell  3' Gene Expression Library Preparation v3.1  with UDI; 10x Genomics Chromium GEM-X Single Cell 3' v4 Gene Expression Library Preparation with UDI; 10x Genomics Chromium GEM-X Single Cell 3' v4 Gene Expression Library Preparation; 10x Genomics Chromium Next GEM Single Cell 3' v3.1 Gene Expression Library Preparation"

# Add the new 10x Genomics Chromium Next GEM v3.1 kit entry to the list.
$ws.Range("B10").Value = $new10xKit
$ws.Range("B11").Value = $new10xKit

# Replace every "--config .../Test/autoAnalysis.multiqc.config.yaml" cell
# with the new path that no longer goes through the Test dir.
$ws.Range("D6").Value = $newConfigPath
$ws.Range("D7").Value = $newConfigPath
$ws.Range("D8").Value = $newConfigPath
$ws.Range("D13").Value = $newConfigPath
$ws.Range("D14").Value = $newConfigPath
$ws.Range("D16").Value = $newConfigPath
$ws.Range("D17").Value = $newConfigPath

# Row height adjustments.
$ws.Rows.Item(4).RowHeight = 50
$ws.Rows.Item(10).RowHeight = 153
$ws.Rows.Item(11).RowHeight = 153

# Update the frozen-pane top-left cell and the active selection.
$win = $excel.ActiveWindow
$panes = $win.Panes
$bottomLeftPane = $panes.Item(2)
$bottomLeftPane.ScrollRow = 5
$bottomLeftPane.ScrollColumn = 1
$ws.Range("C6").Select()
